$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.566.92"
$ws.Range("E2").Value = "  +2.49%  "

$ws.Range("D3").Value = "1.850.90"
$ws.Range("E3").Value = "  +2.04%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.034"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +2.88%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "321.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.29%  "

$ws.Range("E6").Value = "  +2.55%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4387"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.26%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3773"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.26%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07405"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.32%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8744"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.49%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.47"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.95%  "

$ws.Range("D12").Value = "1.854.14"
$ws.Range("E12").Value = "  -8.02%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.526"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.30%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.685"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.63%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07222"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.68%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "82.80"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.62%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.036"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.11%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009031"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.11%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.029"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.46%  "

$ws.Range("E20").Value = "  +1.39%  "

$ws.Range("D21").Value = "27.578.29"
$ws.Range("E21").Value = "  +2.32%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.250"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.30%  "

$ws.Range("E23").Value = "  +2.54%  "

$ws.Range("D24").Value = "2.073.47"
$ws.Range("E24").Value = "  -7.28%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.74"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.61%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.922"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.07%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.69"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.47%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "116.91"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.63%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09047"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.11%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7608"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.35%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.193"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.56%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.497"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.63%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.888"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.19%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.030"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.95%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.147"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.67%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01972"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.60%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05295"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.54%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5146"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.38%  "

$ws.Range("E41").Value = "  +3.05%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1672"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.76%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.731"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.58%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.478"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.87%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "108.72"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.89%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.50"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.88%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.706"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.97%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06403"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.76%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4642"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.46%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.855"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.95%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "39.15"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.06%  "

# Row 28 and 29: coin/link swap, with new price/volume values
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.973"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.42%  "

$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.259"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.91%  "
